$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4575
$ws.Range("I40").Value = 3920
$ws.Range("J40").Value = 5666.6665
$ws.Range("K40").Value = 3920
$ws.Range("L40").Value = 5666.6665
$ws.Range("M40").Value = -3745
$ws.Range("N40").Value = -6016.6665
$ws.Range("H92").Value = 2536.5
$ws.Range("I92").Value = 1069.3334
$ws.Range("J92").Value = 4003.6667
$ws.Range("K92").Value = 1069.3334
$ws.Range("L92").Value = 4003.6667
$ws.Range("M92").Value = 178.6666
$ws.Range("N92").Value = -6499.6667
$ws.Range("H135").Value = 33482.13
$ws.Range("I135").Value = 41112.28
$ws.Range("J135").Value = 1689.8334
$ws.Range("K135").Value = 370010.52
$ws.Range("L135").Value = 15208.5006
$ws.Range("M135").Value = -367475.52
$ws.Range("N135").Value = -20278.5006
$ws.Range("H141").Value = 999
$ws.Range("I141").Value = 999
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 2997
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2183
$ws.Range("N141").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3950
$ws.Range("I2").Value = 5000
$ws.Range("J2").Value = 2900
$ws.Range("K2").Value = 5000
$ws.Range("L2").Value = 2900
$ws.Range("M2").Value = -4887
$ws.Range("N2").Value = -3126
$ws.Range("H61").Value = 38539704
$ws.Range("I61").Value = 47667490
$ws.Range("J61").Value = 203000
$ws.Range("K61").Value = 47667490
$ws.Range("L61").Value = 203000
$ws.Range("M61").Value = -47667278
$ws.Range("N61").Value = -203424
$ws.Range("H92").Value = 19600
$ws.Range("J92").Value = 19600
$ws.Range("L92").Value = 19600
$ws.Range("N92").Value = -24592
$ws.Range("H116").Value = 3950
$ws.Range("I116").Value = 5000
$ws.Range("J116").Value = 2900
$ws.Range("K116").Value = 5000
$ws.Range("L116").Value = 2900
$ws.Range("M116").Value = -2706
$ws.Range("N116").Value = -7488
$ws.Range("H132").Value = 38959.527
$ws.Range("I132").Value = 27833.621
$ws.Range("J132").Value = 64688.188
$ws.Range("K132").Value = 83500.863
$ws.Range("L132").Value = 194064.564
$ws.Range("M132").Value = -80970.863
$ws.Range("N132").Value = -199124.564
$ws.Range("H133").Value = 34300
$ws.Range("J133").Value = 34300
$ws.Range("L133").Value = 34300
$ws.Range("N133").Value = -39360
$ws.Range("H136").Value = 38539704
$ws.Range("I136").Value = 47667490
$ws.Range("J136").Value = 203000
$ws.Range("K136").Value = 143002470
$ws.Range("L136").Value = 609000
$ws.Range("M136").Value = -142999920
$ws.Range("N136").Value = -614100

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3950
$ws.Range("I3").Value = 5000
$ws.Range("J3").Value = 2900
$ws.Range("K3").Value = 5000
$ws.Range("L3").Value = 2900
$ws.Range("M3").Value = -4886
$ws.Range("N3").Value = -3128
$ws.Range("H94").Value = 974.6667
$ws.Range("I94").Value = 967.4286
$ws.Range("K94").Value = 967.4286
$ws.Range("M94").Value = -516.4286
$ws.Range("H134").Value = 2924.244
$ws.Range("I134").Value = 2309.2812
$ws.Range("K134").Value = 6927.8436
$ws.Range("M134").Value = -4392.8436

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 142857810
$ws.Range("I22").Value = 200000540
$ws.Range("K22").Value = 200000540
$ws.Range("M22").Value = -200000190
$ws.Range("H134").Value = 37696.406
$ws.Range("I134").Value = 2657.7083
$ws.Range("J134").Value = 142812.5
$ws.Range("K134").Value = 7973.124899999999
$ws.Range("L134").Value = 428437.5
$ws.Range("M134").Value = -5438.124899999999
$ws.Range("N134").Value = -433507.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 703.29034
$ws.Range("I5").Value = 333.7
$ws.Range("J5").Value = 879.2857
$ws.Range("K5").Value = 1001.1
$ws.Range("L5").Value = 2637.8571
$ws.Range("M5").Value = -889.0999999999999
$ws.Range("N5").Value = -2861.8571
$ws.Range("H113").Value = 596.54346
$ws.Range("I113").Value = 500
$ws.Range("J113").Value = 613.87177
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 1841.61531
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -6181.61531
$ws.Range("H132").Value = 881
$ws.Range("I132").Value = 684.5833
$ws.Range("K132").Value = 6161.2497
$ws.Range("M132").Value = -3631.2497
$ws.Range("H135").Value = 703.29034
$ws.Range("I135").Value = 333.7
$ws.Range("J135").Value = 879.2857
$ws.Range("K135").Value = 3003.3
$ws.Range("L135").Value = 7913.571300000001
$ws.Range("M135").Value = -468.2999999999997
$ws.Range("N135").Value = -12983.5713

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2202.5
$ws.Range("I97").Value = 2202.5
$ws.Range("K97").Value = 2202.5
$ws.Range("M97").Value = -1706.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2201.1428
$ws.Range("I7").Value = 2171.2
$ws.Range("J7").Value = 2800
$ws.Range("K7").Value = 2171.2
$ws.Range("L7").Value = 2800
$ws.Range("M7").Value = -2059.2
$ws.Range("N7").Value = -3024
$ws.Range("H122").Value = 3430.4
$ws.Range("I122").Value = 3151.3333
$ws.Range("J122").Value = 3550
$ws.Range("K122").Value = 9453.999899999999
$ws.Range("L122").Value = 10650
$ws.Range("M122").Value = -7003.999899999999
$ws.Range("N122").Value = -15550
$ws.Range("H126").Value = 2201.1428
$ws.Range("I126").Value = 2171.2
$ws.Range("J126").Value = 2800
$ws.Range("K126").Value = 6513.599999999999
$ws.Range("L126").Value = 8400
$ws.Range("M126").Value = -4043.599999999999
$ws.Range("N126").Value = -13340
$ws.Range("H136").Value = 63375.79
$ws.Range("I136").Value = 35824.17
$ws.Range("J136").Value = 263125
$ws.Range("K136").Value = 107472.51
$ws.Range("L136").Value = 789375
$ws.Range("M136").Value = -104922.51
$ws.Range("N136").Value = -794475

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 9349.714
$ws.Range("I96").Value = 6414.5
$ws.Range("J96").Value = 13263.333
$ws.Range("K96").Value = 6414.5
$ws.Range("L96").Value = 13263.333
$ws.Range("M96").Value = -5041.5
$ws.Range("N96").Value = -16009.333
